$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data rows (2-15 and 17; row 16 is left untouched) get shuffled
# around: the entire content of each source row moves, as a whole, to a new
# row number. Map is sourceRow -> destinationRow.
$map = @{
    2  = 6
    3  = 10
    4  = 11
    5  = 14
    6  = 17
    7  = 12
    8  = 5
    9  = 7
    10 = 8
    11 = 3
    12 = 9
    13 = 15
    14 = 2
    15 = 4
    17 = 13
}

# Only these columns ever actually differ between the rows being permuted
# (the rest -- Lan/Kommun/Provins/Forsamling/Startdatum/.../Rapportor/...
# -- are identical across every row in this sheet, so copying them is a
# no-op; skipping them avoids Excel's automatic type coercion, e.g. turning
# the literal text "2023-08-02" into a real date serial).
$cols = @("A","B","E","F","G","H","K","L","M","N","Q","R","AC","AO")

# Snapshot every affected row's relevant cell values before writing
# anything, since this is a true permutation (every source row is also a
# destination row) and naive in-place copying would clobber data that is
# still needed as a source later on.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $rowData = @{}
    foreach ($col in $cols) {
        $cell = $ws.Range($col + $row)
        $rowData[$col] = $cell.Value()
    }
    $snapshot[$row] = $rowData
}

# Now write each snapshot into its destination row. Cells with no value in
# the source row are cleared (set to empty string) in the destination row.
foreach ($row in $map.Keys) {
    $destRow = $map[$row]
    $rowData = $snapshot[$row]
    foreach ($col in $cols) {
        $val = $rowData[$col]
        if ($null -eq $val) {
            $ws.Range($col + $destRow).Value = ""
        } else {
            $ws.Range($col + $destRow).Value = $val
        }
    }
}
